# Update crypto price/volume columns (D, E) for rows 2-51
# per the latest scrape, matching the commit's refreshed figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.110.61"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.667.47"
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("D5").Value = "'210.61"
$ws.Range("E5").Value = "  -3.35%  "
$ws.Range("D6").Value = "'0.5247"
$ws.Range("E6").Value = "  -2.47%  "
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("D8").Value = "'0.2624"
$ws.Range("E8").Value = "  -3.92%  "
$ws.Range("D9").Value = "'0.06291"
$ws.Range("E9").Value = "  -2.31%  "
$ws.Range("D10").Value = "'21.14"
$ws.Range("E10").Value = "  -2.33%  "
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("D12").Value = "1.670.30"
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("D14").Value = "'0.5539"
$ws.Range("E14").Value = "  -4.31%  "
$ws.Range("D15").Value = "'66.71"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "'0.000007924"
$ws.Range("E16").Value = "  -5.38%  "
$ws.Range("D17").Value = "26.152.39"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "'4.738"
$ws.Range("E19").Value = "  -3.47%  "
$ws.Range("D20").Value = "'186.26"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("E21").Value = "  -4.87%  "
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").Value = "'149.87"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").Value = "'0.1246"
$ws.Range("E25").Value = "  -3.06%  "
$ws.Range("D26").Value = "'7.485"
$ws.Range("E26").Value = "  -4.83%  "
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").Value = "'0.06271"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "'1.354"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("E30").Value = "  -3.67%  "
$ws.Range("D31").Value = "'3.512"
$ws.Range("E31").Value = "  -2.49%  "
$ws.Range("D32").Value = "'3.413"
$ws.Range("D33").Value = "'1.629"
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("D34").Value = "'0.9977"
$ws.Range("E34").Value = "  -3.37%  "
$ws.Range("E35").Value = "  -1.75%  "
$ws.Range("D36").Value = "'2.415"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").Value = "'2.735"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("D38").Value = "'6.118"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "1.106.78"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").Value = "'0.01615"
$ws.Range("E40").Value = "  -2.33%  "
$ws.Range("D41").Value = "'0.8707"
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("D42").Value = "'1.003"
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("D43").Value = "'100.00"
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("D44").Value = "1.820.49"
$ws.Range("E44").Value = "  -1.20%  "
$ws.Range("D45").Value = "'0.00000000111"
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").Value = "'8.010"
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("D49").Value = "'0.05238"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").Value = "'0.4245"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").Value = "'5.965"
$ws.Range("E51").Value = "  -1.28%  "
